$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Row 22 - SUFRE DISCAPACIDAD - "DESCRIPCIÓN" column (C): fix "2: No" -> "0: No"
$ws.Range("C22").Value = "Indicador si el Alumno sufre o no discapacidad. Valores:" + $nl + "1: Si" + $nl + "0: No" + $nl + "Vacío: No"

# Row 25 - ¿CUENTA CON ESCALA DE PAGO? - "DESCRIPCIÓN" column (C): fix "2: No" -> "0: No"
$ws.Range("C25").Value = "Indicador si el Alumno cuenta con una Escala de Pago. Valores:" + $nl + "1: Si" + $nl + "0: No" + $nl + "Vacío: No"

# Row 30 - ¿REALIZÓ MOVILIDAD NACIONAL? - "DESCRIPCIÓN" column (C): fix "2: No" -> "0: No"
$ws.Range("C30").Value = "Indicador si el alumno realizó o no Movilidad Nacional. Valores:" + $nl + "1: Si" + $nl + "0: No" + $nl + "Vacío: No"

# Row 34 - ¿REALIZÓ MOVILIDAD INTERNACIONAL? - "DESCRIPCIÓN" column (C): fix "2: No" -> "0: No"
$ws.Range("C34").Value = "Indicador si el alumno realizó o no Movilidad Internacional. Valores:" + $nl + "1: Si" + $nl + "0: No" + $nl + "Vacío: No"

# Row 38 - ¿ES ALUMNO EN RIESGO? - "DESCRIPCIÓN" column (C): fix "2: No" -> "0: No"
$ws.Range("C38").Value = "Indicador si el alumno se encuentra en riesgo. Valores:" + $nl + "1: Si" + $nl + "0: No" + $nl + "Vacío: No"
